# Insert a new week of "Tuna" price records (date 2022-05-26 / serial 44707)
# right before the existing row 18, pushing the remaining records down by
# three rows (dimension grows from T95 to T98).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows at position 18; everything below shifts down.
$ws.Rows("18:20").Insert()

# Columns that are constant across every data row in this sheet.
$ws.Range("A18:A20").Value = 2
$ws.Range("B18:B20").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C18:C20").Value = "Coquimbo"
$ws.Range("D18:D20").Value = 44707
$ws.Range("D18:D20").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E18:E20").Value = 4
$ws.Range("F18:F20").Value = "Fruta"
$ws.Range("G18:G20").Value = 100107
$ws.Range("H18:H20").Value = "Otros"
$ws.Range("I18:I20").Value = 100107011
$ws.Range("J18:J20").Value = "Tuna"
$ws.Range("K18:K20").Value = "Sin especificar"
$ws.Range("Q18:Q20").Value = "$/caja 18 kilos"
$ws.Range("R18:R20").Value = "Provincia de Limarí"
$ws.Range("T18:T20").Value = 18

# Row 18: Especial
$ws.Range("L18").Value = "Especial"
$ws.Range("M18").Value = 140
$ws.Range("N18").Value = 16000
$ws.Range("O18").Value = 17000
$ws.Range("P18").Value = 16500
$ws.Range("S18").Value = 917

# Row 19: Primera
$ws.Range("L19").Value = "Primera"
$ws.Range("M19").Value = 300
$ws.Range("N19").Value = 12000
$ws.Range("O19").Value = 13000
$ws.Range("P19").Value = 12500
$ws.Range("S19").Value = 694

# Row 20: Segunda
$ws.Range("L20").Value = "Segunda"
$ws.Range("M20").Value = 100
$ws.Range("N20").Value = 8000
$ws.Range("O20").Value = 9000
$ws.Range("P20").Value = 8500
$ws.Range("S20").Value = 472
